$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.358839897227597
$ws.Range("C2").Value = 0.298285615648922
$ws.Range("E2").Value = 0.05035469001898107
$ws.Range("F2").Value = 3.706044695401189
$ws.Range("G2").Value = 0.002611361463643396
$ws.Range("J2").Value = 0.2062711819693774
$ws.Range("L2").Value = 0.239546571594353
$ws.Range("M2").Value = 0.5867500080686696
$ws.Range("N2").Value = 2.912624198134026
$ws.Range("B3").Value = 3.24715315543159
$ws.Range("C3").Value = 0.2659263982062612
$ws.Range("E3").Value = 0.04856748507910957
$ws.Range("F3").Value = 3.67981311465897
$ws.Range("G3").Value = 0.002617772876892005
$ws.Range("J3").Value = 0.2073141481491483
$ws.Range("L3").Value = 0.239272721205019
$ws.Range("M3").Value = 0.5732557173137991
$ws.Range("N3").Value = 2.923700254840469
$ws.Range("B4").Value = 3.180698909940588
$ws.Range("C4").Value = 0.2461972570567639
$ws.Range("E4").Value = 0.04745524510442145
$ws.Range("F4").Value = 3.665464287824889
$ws.Range("G4").Value = 0.002621915001888303
$ws.Range("J4").Value = 0.2080127221154555
$ws.Range("L4").Value = 0.2391944304334075
$ws.Range("M4").Value = 0.5653181618475713
$ws.Range("N4").Value = 2.931260879319098
$ws.Range("B5").Value = 3.154149411195249
$ws.Range("C5").Value = 0.2381915310794795
$ws.Range("E5").Value = 0.04699815885523684
$ws.Range("F5").Value = 3.66005777879576
$ws.Range("G5").Value = 0.002623654809990938
$ws.Range("J5").Value = 0.2083119784844882
$ws.Range("L5").Value = 0.2391851987342619
$ws.Range("M5").Value = 0.5621708611408209
$ws.Range("N5").Value = 2.934532333610946
$ws.Range("B6").Value = 3.149772897634421
$ws.Range("C6").Value = 0.2368642135626828
$ws.Range("E6").Value = 0.04692202536404189
$ws.Range("F6").Value = 3.659186618671058
$ws.Range("G6").Value = 0.002623946841324957
$ws.Range("J6").Value = 0.2083625490316336
$ws.Range("L6").Value = 0.239185037482244
$ws.Range("M6").Value = 0.5616535226832227
$ws.Range("N6").Value = 2.935087039506683
$ws.Range("B7").Value = 3.180338706662781
$ws.Range("C7").Value = 0.2460891524876274
$ws.Range("E7").Value = 0.04744909635097372
$ws.Range("F7").Value = 3.665389590678075
$ws.Range("G7").Value = 0.002621938255440698
$ws.Range("J7").Value = 0.2080166990150154
$ws.Range("L7").Value = 0.239194214027016
$ws.Range("M7").Value = 0.5652753629447034
$ws.Range("N7").Value = 2.931304229051719
$ws.Range("B8").Value = 3.319888403889763
$ws.Range("C8").Value = 0.2870985603576628
$ws.Range("E8").Value = 0.04974149387580162
$ws.Range("F8").Value = 3.69663455315154
$ws.Range("G8").Value = 0.002613529581121085
$ws.Range("J8").Value = 0.206618689231064
$ws.Range("L8").Value = 0.2394335444672109
$ws.Range("M8").Value = 0.5820248293571524
$ws.Range("N8").Value = 2.916285122020398
$ws.Range("B9").Value = 3.610510425707446
$ws.Range("C9").Value = 0.3686760096543367
$ws.Range("E9").Value = 0.0541232167500354
$ws.Range("F9").Value = 3.771911695172832
$ws.Range("G9").Value = 0.002598662184665781
$ws.Range("J9").Value = 0.2043411968546884
$ws.Range("L9").Value = 0.2406126947912881
$ws.Range("M9").Value = 0.617642964631159
$ws.Range("N9").Value = 2.892891337760844
$ws.Range("B10").Value = 3.83458199780307
$ws.Range("C10").Value = 0.4293936843794768
$ws.Range("E10").Value = 0.05727961681637339
$ws.Range("F10").Value = 3.83585255966716
$ws.Range("G10").Value = 0.002588716001425442
$ws.Range("J10").Value = 0.2029540338263764
$ws.Range("L10").Value = 0.2419079590908595
$ws.Range("M10").Value = 0.6455215318188792
$ws.Range("N10").Value = 2.879438730277798
$ws.Range("B11").Value = 3.93885699771414
$ws.Range("C11").Value = 0.4572028139264717
$ws.Range("E11").Value = 0.05870331529216699
$ws.Range("F11").Value = 3.866837839236638
$ws.Range("G11").Value = 0.002584400799340494
$ws.Range("J11").Value = 0.2023858291827523
$ws.Range("L11").Value = 0.2425895918112388
$ws.Range("M11").Value = 0.6585801503109323
$ws.Range("N11").Value = 2.874138883763138
$ws.Range("B12").Value = 3.978683676243406
$ws.Range("C12").Value = 0.4677617813760548
$ws.Range("E12").Value = 0.05924079586505471
$ws.Range("F12").Value = 3.878845751250566
$ws.Range("G12").Value = 0.002582796655105787
$ws.Range("J12").Value = 0.2021797611329248
$ws.Range("L12").Value = 0.2428609288759276
$ws.Range("M12").Value = 0.6635795666180897
$ws.Range("N12").Value = 2.872250621417365
$ws.Range("B13").Value = 3.970091118679363
$ws.Range("C13").Value = 0.4654864434195929
$ws.Range("E13").Value = 0.05912511118330954
$ws.Range("F13").Value = 3.876247401782194
$ws.Range("G13").Value = 0.002583140808221842
$ws.Range("J13").Value = 0.202223735856073
$ws.Range("L13").Value = 0.2428019048480934
$ws.Range("M13").Value = 0.6625004291586691
$ws.Range("N13").Value = 2.872652001755867
$ws.Range("B14").Value = 3.942126731124347
$ws.Range("C14").Value = 0.458070933528461
$ws.Range("E14").Value = 0.05874756646092294
$ws.Range("F14").Value = 3.867820228153278
$ws.Range("G14").Value = 0.002584268226341108
$ws.Range("J14").Value = 0.2023686931111399
$ws.Range("L14").Value = 0.242611650377647
$ws.Range("M14").Value = 0.6589903635122738
$ws.Range("N14").Value = 2.873981152324774
$ws.Range("B15").Value = 3.925042103817759
$ws.Range("C15").Value = 0.4535324340601505
$ws.Range("E15").Value = 0.05851609860259543
$ws.Range("F15").Value = 3.862694128495036
$ws.Range("G15").Value = 0.002584962695215461
$ws.Range("J15").Value = 0.2024586704882694
$ws.Range("L15").Value = 0.2424968332499589
$ws.Range("M15").Value = 0.656847439376989
$ws.Range("N15").Value = 2.874810773889337
$ws.Range("B16").Value = 3.827814795531765
$ws.Range("C16").Value = 0.4275801819960066
$ws.Range("E16").Value = 0.05718633903433457
$ws.Range("F16").Value = 3.833865923359667
$ws.Range("G16").Value = 0.002589002207851628
$ws.Range("J16").Value = 0.2029924370017184
$ws.Range("L16").Value = 0.2418652662312653
$ws.Range("M16").Value = 0.6446757187380285
$ws.Range("N16").Value = 2.879801651390423
$ws.Range("B17").Value = 3.768771366529677
$ws.Range("C17").Value = 0.4117083835801623
$ws.Range("E17").Value = 0.05636753414101392
$ws.Range("F17").Value = 3.816667993821227
$ws.Range("G17").Value = 0.002591533815716908
$ws.Range("J17").Value = 0.2033360215628832
$ws.Range("L17").Value = 0.2415014379955025
$ws.Range("M17").Value = 0.6373053783038856
$ws.Range("N17").Value = 2.883073922504721
$ws.Range("B18").Value = 3.735031567129852
$ws.Range("C18").Value = 0.4025970027590233
$ws.Range("E18").Value = 0.05589542788714752
$ws.Range("F18").Value = 3.80695478161087
$ws.Range("G18").Value = 0.002593009646627378
$ws.Range("J18").Value = 0.2035395515201976
$ws.Range("L18").Value = 0.2413008759003858
$ws.Range("M18").Value = 0.6331015771125337
$ws.Range("N18").Value = 2.885033163183692
$ws.Range("B19").Value = 3.723645618571084
$ws.Range("C19").Value = 0.3995150432281207
$ws.Range("E19").Value = 0.05573537987361377
$ws.Range("F19").Value = 3.803696681547791
$ws.Range("G19").Value = 0.002593512729597292
$ws.Range("J19").Value = 0.2036094762253704
$ws.Range("L19").Value = 0.2412344658283061
$ws.Range("M19").Value = 0.6316843186054371
$ws.Range("N19").Value = 2.885709749515811
$ws.Range("B20").Value = 3.775033806472777
$ws.Range("C20").Value = 0.4133961276929767
$ws.Range("E20").Value = 0.05645481571529132
$ws.Range("F20").Value = 3.81848024952879
$ws.Range("G20").Value = 0.002591262282409175
$ws.Range("J20").Value = 0.2032988343547366
$ws.Range("L20").Value = 0.2415392679971333
$ws.Range("M20").Value = 0.6380862955706377
$ws.Range("N20").Value = 2.882717596697717
$ws.Range("B21").Value = 3.950331298045285
$ws.Range("C21").Value = 0.4602482723377648
$ws.Range("E21").Value = 0.05885850424011352
$ws.Range("F21").Value = 3.870288032526389
$ws.Range("G21").Value = 0.002583936265541774
$ws.Range("J21").Value = 0.2023258681902824
$ws.Range("L21").Value = 0.2426671745667761
$ws.Range("M21").Value = 0.6600198763489118
$ws.Range("N21").Value = 2.873587521462184
$ws.Range("B22").Value = 4.066881553000087
$ws.Range("C22").Value = 0.4910340969252616
$ws.Range("E22").Value = 0.06041992640053095
$ws.Range("F22").Value = 3.905747887790596
$ws.Range("G22").Value = 0.002579322660851886
$ws.Range("J22").Value = 0.2017430440576931
$ws.Range("L22").Value = 0.2434813365990394
$ws.Range("M22").Value = 0.6746719560270478
$ws.Range("N22").Value = 2.868312588811875
$ws.Range("B23").Value = 4.004494021504968
$ws.Range("C23").Value = 0.4745876125028303
$ws.Range("E23").Value = 0.05958740234840576
$ws.Range("F23").Value = 3.886675350462326
$ws.Range("G23").Value = 0.002581769131895388
$ws.Range("J23").Value = 0.2020492309268249
$ws.Range("L23").Value = 0.2430397797311059
$ws.Range("M23").Value = 0.6668227522247534
$ws.Range("N23").Value = 2.87106432951245
$ws.Range("B24").Value = 3.772201918899555
$ws.Range("C24").Value = 0.4126330565323997
$ws.Range("E24").Value = 0.05641535997327551
$ws.Range("F24").Value = 3.817660386532424
$ws.Range("G24").Value = 0.002591384978972516
$ws.Range("J24").Value = 0.2033156280116195
$ws.Range("L24").Value = 0.241522138237066
$ws.Range("M24").Value = 0.6377331385308054
$ws.Range("N24").Value = 2.882878448800156
$ws.Range("B25").Value = 3.530051854809074
$ws.Range("C25").Value = 0.3464752972495262
$ws.Range("E25").Value = 0.05294935202497797
$ws.Range("F25").Value = 3.750037958701796
$ws.Range("G25").Value = 0.002602511798394809
$ws.Range("J25").Value = 0.204907314400117
$ws.Range("L25").Value = 0.2402181074471414
$ws.Range("M25").Value = 0.6077084942043101
$ws.Range("N25").Value = 2.898567128365741

Write-Output "Updated 216 cells for Case_2_0 (380 kV)"